$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "26.490.87"
$ws.Cells.Item(2, 5).Value = "  +0.51%  "

$ws.Cells.Item(3, 4).Value = "1.836.04"
$ws.Cells.Item(3, 5).Value = "  +0.03%  "

$ws.Cells.Item(4, 5).Value = "  -0.11%  "

$ws.Cells.Item(5, 4).Value = "'259.38"
$ws.Cells.Item(5, 5).Value = "  +0.15%  "

$ws.Cells.Item(6, 5).Value = "  -0.05%  "

$ws.Cells.Item(7, 4).Value = "'0.5263"
$ws.Cells.Item(7, 5).Value = "  +1.42%  "

$ws.Cells.Item(8, 4).Value = "'0.3196"
$ws.Cells.Item(8, 5).Value = "  -1.30%  "

$ws.Cells.Item(9, 4).Value = "'0.06781"
$ws.Cells.Item(9, 5).Value = "  +0.49%  "

$ws.Cells.Item(10, 4).Value = "'18.72"
$ws.Cells.Item(10, 5).Value = "  +0.93%  "

$ws.Cells.Item(11, 4).Value = "'0.7806"
$ws.Cells.Item(11, 5).Value = "  +2.39%  "

$ws.Cells.Item(12, 4).Value = "'0.07729"
$ws.Cells.Item(12, 5).Value = "  +1.07%  "

$ws.Cells.Item(13, 4).Value = "1.840.14"
$ws.Cells.Item(13, 5).Value = "  +0.81%  "

$ws.Cells.Item(14, 5).Value = "  -1.07%  "

$ws.Cells.Item(15, 4).Value = "'5.005"
$ws.Cells.Item(15, 5).Value = "  -0.12%  "

$ws.Cells.Item(16, 5).Value = "  -0.25%  "

$ws.Cells.Item(17, 4).Value = "'13.81"
$ws.Cells.Item(17, 5).Value = "  -0.99%  "

$ws.Cells.Item(18, 4).Value = "'1.001"
$ws.Cells.Item(18, 5).Value = "  -0.02%  "

$ws.Cells.Item(19, 4).Value = "'0.000007928"
$ws.Cells.Item(19, 5).Value = "  +0.77%  "

$ws.Cells.Item(20, 4).Value = "26.510.75"
$ws.Cells.Item(20, 5).Value = "  +0.38%  "

$ws.Cells.Item(21, 4).Value = "2.071.56"
$ws.Cells.Item(21, 5).Value = "  -0.16%  "

$ws.Cells.Item(22, 4).Value = "'4.603"
$ws.Cells.Item(22, 5).Value = "  +1.18%  "

$ws.Cells.Item(23, 2).Value = "Chainlink"
$ws.Cells.Item(23, 3).Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Cells.Item(23, 4).Value = "'5.964"
$ws.Cells.Item(23, 5).Value = "  +0.62%  "

$ws.Cells.Item(24, 2).Value = "Cosmos"
$ws.Cells.Item(24, 3).Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Cells.Item(24, 4).Value = "'9.362"
$ws.Cells.Item(24, 5).Value = "  -0.48%  "

$ws.Cells.Item(25, 4).Value = "'141.56"
$ws.Cells.Item(25, 5).Value = "  -1.79%  "

$ws.Cells.Item(26, 4).Value = "'2.146"
$ws.Cells.Item(26, 5).Value = "  -4.05%  "

$ws.Cells.Item(27, 4).Value = "'1.679"
$ws.Cells.Item(27, 5).Value = "  +2.19%  "

$ws.Cells.Item(28, 5).Value = "  +0.04%  "

$ws.Cells.Item(29, 4).Value = "'111.57"
$ws.Cells.Item(29, 5).Value = "  +0.28%  "

$ws.Cells.Item(30, 4).Value = "'4.129"
$ws.Cells.Item(30, 5).Value = "  -0.84%  "

$ws.Cells.Item(31, 4).Value = "'0.08684"
$ws.Cells.Item(31, 5).Value = "  -0.19%  "

$ws.Cells.Item(32, 4).Value = "'4.057"
$ws.Cells.Item(32, 5).Value = "  -2.04%  "

$ws.Cells.Item(33, 4).Value = "'0.04845"
$ws.Cells.Item(33, 5).Value = "  +1.34%  "

$ws.Cells.Item(34, 4).Value = "'0.7279"
$ws.Cells.Item(34, 5).Value = "  +5.69%  "

$ws.Cells.Item(35, 4).Value = "'1.128"
$ws.Cells.Item(35, 5).Value = "  +0.74%  "

$ws.Cells.Item(36, 4).Value = "'2.852"
$ws.Cells.Item(36, 5).Value = "  +0.03%  "

$ws.Cells.Item(37, 4).Value = "'3.086"
$ws.Cells.Item(37, 5).Value = "  +0.88%  "

$ws.Cells.Item(38, 4).Value = "'2.235"
$ws.Cells.Item(38, 5).Value = "  +2.06%  "

$ws.Cells.Item(39, 4).Value = "'0.01762"
$ws.Cells.Item(39, 5).Value = "  +0.25%  "

$ws.Cells.Item(40, 4).Value = "'0.4743"
$ws.Cells.Item(40, 5).Value = "  -1.49%  "

$ws.Cells.Item(41, 4).Value = "'0.8909"
$ws.Cells.Item(41, 5).Value = "  +0.35%  "

$ws.Cells.Item(42, 4).Value = "'109.41"
$ws.Cells.Item(42, 5).Value = "  -1.21%  "

$ws.Cells.Item(43, 4).Value = "'5.903"
$ws.Cells.Item(43, 5).Value = "  -3.01%  "

$ws.Cells.Item(44, 5).Value = "  +0.03%  "

$ws.Cells.Item(45, 4).Value = "'7.621"
$ws.Cells.Item(45, 5).Value = "  -0.39%  "

$ws.Cells.Item(46, 4).Value = "'0.4126"
$ws.Cells.Item(46, 5).Value = "  -0.03%  "

$ws.Cells.Item(47, 4).Value = "'0.05851"
$ws.Cells.Item(47, 5).Value = "  -0.09%  "

$ws.Cells.Item(48, 4).Value = "'8.923"
$ws.Cells.Item(48, 5).Value = "  -0.75%  "

$ws.Cells.Item(49, 2).Value = "Algorand"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Cells.Item(49, 4).Value = "'0.1225"
$ws.Cells.Item(49, 5).Value = "  -0.50%  "

$ws.Cells.Item(50, 2).Value = "Elrond"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Cells.Item(50, 4).Value = "'34.72"
$ws.Cells.Item(50, 5).Value = "  +0.12%  "

$ws.Cells.Item(51, 4).Value = "'0.8925"
$ws.Cells.Item(51, 5).Value = "  +1.07%  "
